$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 6, pushing the old 6-cylinder / 8-cylinder / footer
# rows down by one. This mirrors the existing 6-cylinder group (which already
# spans two merged rows) so that the 4-cylinder group also gets a top/bottom
# pair of rows (split further by transmission type). Excel's row insert
# already shifts the old A6:A7 merge down to A7:A8, so that one needs no
# further work - only the new row 5/6 pair needs to be merged.
$ws.Rows("6").Insert()
$ws.Range("A5:A6").Merge()

# Fix up column A's formatting for the new pairing: row 5 becomes the top of
# the merged pair (vertical-top alignment, like the existing 6-cylinder
# group's top cell) and row 6 becomes the blank bottom half (same plain
# border style as the existing 6-cylinder group's bottom cell). Copying the
# formats from the equivalent cells below keeps the style table unchanged
# instead of growing it with near-duplicate entries.
$ws.Range("A7").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A8").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$ws.Range("C7").Copy()
$ws.Range("C6").PasteSpecial(-4122)
$ws.Range("H7").Copy()
$ws.Range("H6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 5: 4-cylinder, automatic transmission (top half of the new merge) ---
$ws.Range("B5").Value = 0
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 91
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = 2.14
$ws.Range("G5").Value = ""

# --- Row 6: 4-cylinder, manual transmission (bottom half of the new merge) ---
$ws.Range("A6").Value = ""
$ws.Range("B6").Value = 1
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = 81.8
$ws.Range("E6").Value = 21.87235698318771
$ws.Range("F6").Value = 2.3003
$ws.Range("G6").Value = 0.5982073312080948

# --- Row 7: 6-cylinder, automatic transmission (top half of existing merge) ---
$ws.Range("B7").Value = 0
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = 131.6666666666667
$ws.Range("E7").Value = 37.52776749732568
$ws.Range("F7").Value = 2.755
$ws.Range("G7").Value = 0.1281600561797629

# --- Row 8: 6-cylinder, manual transmission (bottom half of existing merge) ---
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = 115.25
$ws.Range("E8").Value = 9.178779875342908
$ws.Range("F8").Value = 3.38875
$ws.Range("G8").Value = 0.1162163929916946

# --- Row 9: 8-cylinder (single row, like before) ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 0
$ws.Range("C9").Value = 14
$ws.Range("D9").Value = 209.2142857142857
$ws.Range("E9").Value = 50.97688551827051
$ws.Range("F9").Value = 3.999214285714287
$ws.Range("G9").Value = 0.7594047444769265
